$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure D (Price) and E (Volume) columns keep their original text formatting
# (Excel would otherwise silently coerce numeric-looking strings to numbers,
# which can drop trailing zeros, e.g. '2.10' -> 2.1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.871.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.61'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.511'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.87%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0626'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.97'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.74'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.650.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.74'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.856.18'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.99'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.79'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.56'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.51%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.255.10'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.833'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.776.07'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.10'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.66'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.52'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.57'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0514'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.64'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0961'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.13%  '
